$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 907.4
$ws.Range("I2").Value = 418.5
$ws.Range("J2").Value = 1233.3334
$ws.Range("K2").Value = 418.5
$ws.Range("L2").Value = 1233.3334
$ws.Range("M2").Value = -305.5
$ws.Range("N2").Value = -1459.3334

$ws.Range("H18").Value = 3837.8
$ws.Range("I18").Value = 3397.25
$ws.Range("K18").Value = 3397.25
$ws.Range("M18").Value = -3113.25

$ws.Range("H70").Value = 2411.25
$ws.Range("I70").Value = 1936.25
$ws.Range("K70").Value = 5808.75
$ws.Range("M70").Value = -5538.75

$ws.Range("H73").Value = 2411.25
$ws.Range("I73").Value = 1936.25
$ws.Range("K73").Value = 5808.75
$ws.Range("M73").Value = -4872.75

$ws.Range("H100").Value = 3247.739
$ws.Range("I100").Value = 1088.2667
$ws.Range("J100").Value = 7296.75
$ws.Range("K100").Value = 1088.2667
$ws.Range("L100").Value = 7296.75
$ws.Range("M100").Value = -547.2666999999999
$ws.Range("N100").Value = -8378.75

$ws.Range("H117").Value = 109240
$ws.Range("J117").Value = 109240
$ws.Range("L117").Value = 109240
$ws.Range("N117").Value = -118418

$ws.Range("H138").Value = 2858.1619
$ws.Range("I138").Value = 2419.375
$ws.Range("J138").Value = 2916.6667
$ws.Range("K138").Value = 7258.125
$ws.Range("L138").Value = 8750.000100000001
$ws.Range("M138").Value = -2118.125
$ws.Range("N138").Value = -19030.0001

$ws.Range("H141").Value = 14578.333
$ws.Range("I141").Value = 13666
$ws.Range("K141").Value = 40998
$ws.Range("M141").Value = -35818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7697989
$ws.Range("I32").Value = 9806270
$ws.Range("K32").Value = 9806270
$ws.Range("M32").Value = -9805983

$ws.Range("H45").Value = 27779970
$ws.Range("I45").Value = 35716110
$ws.Range("K45").Value = 35716110
$ws.Range("M45").Value = -35715733

$ws.Range("H101").Value = 45066.668
$ws.Range("J101").Value = 45066.668
$ws.Range("L101").Value = 45066.668
$ws.Range("N101").Value = -51556.668

$ws.Range("H104").Value = 32333.334
$ws.Range("J104").Value = 32333.334
$ws.Range("L104").Value = 32333.334
$ws.Range("N104").Value = -39321.334

$ws.Range("H106").Value = 53395
$ws.Range("J106").Value = 53395
$ws.Range("L106").Value = 53395
$ws.Range("N106").Value = -55919

$ws.Range("H110").Value = 1534.7858
$ws.Range("I110").Value = 1338.3
$ws.Range("J110").Value = 2026
$ws.Range("K110").Value = 1338.3
$ws.Range("L110").Value = 2026
$ws.Range("M110").Value = 706.7
$ws.Range("N110").Value = -6116

$ws.Range("H122").Value = 2445
$ws.Range("I122").Value = 1726.35
$ws.Range("J122").Value = 3642.75
$ws.Range("K122").Value = 5179.049999999999
$ws.Range("L122").Value = 10928.25
$ws.Range("M122").Value = -2729.049999999999
$ws.Range("N122").Value = -15828.25

$ws.Range("H132").Value = 3311.7188
$ws.Range("I132").Value = 2417.6667
$ws.Range("K132").Value = 7253.000100000001
$ws.Range("M132").Value = -4723.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 99
$ws.Range("J11").Value = 99
$ws.Range("L11").Value = 99
$ws.Range("N11").Value = -379

$ws.Range("H134").Value = 335307.34
$ws.Range("I134").Value = 1434.8572
$ws.Range("J134").Value = 2204993.2
$ws.Range("K134").Value = 4304.571599999999
$ws.Range("L134").Value = 6614979.600000001
$ws.Range("M134").Value = -1769.571599999999
$ws.Range("N134").Value = -6620049.600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 785749.4
$ws.Range("J31").Value = 1063898.8
$ws.Range("L31").Value = 1063898.8
$ws.Range("N31").Value = -1064488.8

$ws.Range("H34").Value = 785749.4
$ws.Range("J34").Value = 1063898.8
$ws.Range("L34").Value = 1063898.8
$ws.Range("N34").Value = -1064302.8

$ws.Range("H62").Value = 837337.7
$ws.Range("I62").Value = 1254662
$ws.Range("K62").Value = 1254662
$ws.Range("M62").Value = -1254038

$ws.Range("H65").Value = 837337.7
$ws.Range("I65").Value = 1254662
$ws.Range("K65").Value = 6273310
$ws.Range("M65").Value = -6270190

$ws.Range("H105").Value = 1535.4445
$ws.Range("I105").Value = 1304.6666
$ws.Range("J105").Value = 1997
$ws.Range("K105").Value = 1304.6666
$ws.Range("L105").Value = 1997
$ws.Range("M105").Value = 442.3334
$ws.Range("N105").Value = -5491

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 18357878
$ws.Range("I4").Value = 11255459
$ws.Range("K4").Value = 33766377
$ws.Range("M4").Value = -33766265

$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -14717
$ws.Range("N32").Value = $null

$ws.Range("H131").Value = 35580.43
$ws.Range("I131").Value = 20010
$ws.Range("J131").Value = 47258.25
$ws.Range("K131").Value = 60030
$ws.Range("L131").Value = 141774.75
$ws.Range("M131").Value = -54990
$ws.Range("N131").Value = -151854.75

$ws.Range("H134").Value = 11643.357
$ws.Range("J134").Value = 12850.56
$ws.Range("L134").Value = 38551.68
$ws.Range("N134").Value = -48691.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3123.8
$ws.Range("I102").Value = 2579.0667
$ws.Range("K102").Value = 2579.0667
$ws.Range("M102").Value = -957.0666999999999

$ws.Range("H105").Value = 48995
$ws.Range("J105").Value = 48995
$ws.Range("L105").Value = 48995
$ws.Range("N105").Value = -55983

$ws.Range("H109").Value = 45213.332
$ws.Range("J109").Value = 45213.332
$ws.Range("L109").Value = 45213.332
$ws.Range("N109").Value = -47293.332

$ws.Range("H113").Value = 3854.65
$ws.Range("I113").Value = 2611.3333
$ws.Range("J113").Value = 4871.909
$ws.Range("K113").Value = 2611.3333
$ws.Range("L113").Value = 4871.909
$ws.Range("M113").Value = -441.3332999999998
$ws.Range("N113").Value = -9211.909

$ws.Range("H126").Value = 4178.357
$ws.Range("I126").Value = 3585.3333
$ws.Range("K126").Value = 10755.9999
$ws.Range("M126").Value = -8285.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1731.591
$ws.Range("I22").Value = 1841.3334
$ws.Range("J22").Value = 1599.9
$ws.Range("K22").Value = 1841.3334
$ws.Range("L22").Value = 1599.9
$ws.Range("M22").Value = -1546.3334
$ws.Range("N22").Value = -2189.9

$ws.Range("H27").Value = 1731.591
$ws.Range("I27").Value = 1841.3334
$ws.Range("J27").Value = 1599.9
$ws.Range("K27").Value = 1841.3334
$ws.Range("L27").Value = 1599.9
$ws.Range("M27").Value = -1734.3334
$ws.Range("N27").Value = -1813.9

$ws.Range("H46").Value = 3362.611
$ws.Range("I46").Value = 2534.3333
$ws.Range("J46").Value = 4190.8887
$ws.Range("K46").Value = 2534.3333
$ws.Range("L46").Value = 4190.8887
$ws.Range("M46").Value = -2346.3333
$ws.Range("N46").Value = -4566.8887

$ws.Range("H61").Value = 2892.0667
$ws.Range("J61").Value = 3999.6667
$ws.Range("L61").Value = 3999.6667
$ws.Range("N61").Value = -4403.6667

$ws.Range("H93").Value = 35715830
$ws.Range("I93").Value = 62501250
$ws.Range("K93").Value = 62501250
$ws.Range("M93").Value = -62500002

$ws.Range("H101").Value = 150181
$ws.Range("J101").Value = 150181
$ws.Range("L101").Value = 150181
$ws.Range("N101").Value = -156671

$ws.Range("H105").Value = 56000
$ws.Range("J105").Value = 56000
$ws.Range("L105").Value = 56000
$ws.Range("N105").Value = -62988

$ws.Range("H113").Value = 2892.0667
$ws.Range("J113").Value = 3999.6667
$ws.Range("L113").Value = 3999.6667
$ws.Range("N113").Value = -8339.6667

$ws.Range("H122").Value = 7093.5947
$ws.Range("I122").Value = 7331.4165
$ws.Range("J122").Value = 6654.5386
$ws.Range("K122").Value = 21994.2495
$ws.Range("L122").Value = 19963.6158
$ws.Range("M122").Value = -19544.2495
$ws.Range("N122").Value = -24863.6158

$ws.Range("H132").Value = 427431.38
$ws.Range("I132").Value = 89044.75
$ws.Range("K132").Value = 267134.25
$ws.Range("M132").Value = -264604.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 77841
$ws.Range("J103").Value = 77841
$ws.Range("L103").Value = 77841
$ws.Range("N103").Value = -80185

$ws.Range("H122").Value = 2472.4707
$ws.Range("I122").Value = 2372.25
$ws.Range("J122").Value = 2713
$ws.Range("K122").Value = 7116.75
$ws.Range("L122").Value = 8139
$ws.Range("M122").Value = -4666.75
$ws.Range("N122").Value = -13039

$ws.Range("H126").Value = 3075
$ws.Range("I126").Value = 2971.6667
$ws.Range("J126").Value = 4005
$ws.Range("K126").Value = 8915.000100000001
$ws.Range("L126").Value = 12015
$ws.Range("M126").Value = -6445.000100000001
$ws.Range("N126").Value = -16955

$ws.Range("H136").Value = 5284.963
$ws.Range("I136").Value = 5558.476
$ws.Range("J136").Value = 4327.6665
$ws.Range("K136").Value = 16675.428
$ws.Range("L136").Value = 12982.9995
$ws.Range("M136").Value = -14125.428
$ws.Range("N136").Value = -18082.9995
